# Workbook: "Hortaliza, Vega Monumental Concepción - Zanahoria"
# Change: insert 2 new weekly records at the top of the existing data block
# (rows 473-474), pushing all subsequent rows down by 2. The dataset keeps
# growing weekly, and this commit adds the newest week's observations while
# preserving every previously recorded row (shifted to rows 475-503).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the first data row of this block (473),
# pushing the old rows 473..501 down to 475..503.
$ws.Rows.Item(473).Insert()
$ws.Rows.Item(473).Insert()

$ws.Range("A473").Value = 11
$ws.Range("B473").Value = "Vega Monumental Concepción"
$ws.Range("C473").Value = "Bíobío"
$ws.Range("D473").Value = 45267
$ws.Range("E473").Value = 8
$ws.Range("F473").Value = 100114013
$ws.Range("G473").Value = "Zanahoria"
$ws.Range("H473").Value = "Sin especificar"
$ws.Range("I473").Value = "Primera"
$ws.Range("J473").Value = 800
$ws.Range("K473").Value = 5000
$ws.Range("L473").Value = 5500
$ws.Range("M473").Value = 5250
$ws.Range("N473").Value = "`$/saco 20 kilos"
$ws.Range("O473").Value = "Región Metropolitana"
$ws.Range("P473").Value = 262
$ws.Range("Q473").Value = 20
$ws.Range("R473").Value = "Hortaliza"
$ws.Range("A474").Value = 11
$ws.Range("B474").Value = "Vega Monumental Concepción"
$ws.Range("C474").Value = "Bíobío"
$ws.Range("D474").Value = 45267
$ws.Range("E474").Value = 8
$ws.Range("F474").Value = 100114013
$ws.Range("G474").Value = "Zanahoria"
$ws.Range("H474").Value = "Sin especificar"
$ws.Range("I474").Value = "Segunda"
$ws.Range("J474").Value = 400
$ws.Range("K474").Value = 4500
$ws.Range("L474").Value = 4500
$ws.Range("M474").Value = 4500
$ws.Range("N474").Value = "`$/saco 20 kilos"
$ws.Range("O474").Value = "Región Metropolitana"
$ws.Range("P474").Value = 225
$ws.Range("Q474").Value = 20
$ws.Range("R474").Value = "Hortaliza"

